$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation is inserted as row 35, pushing the existing
# rows 35-149 down to 36-150 (matching the row-shift pattern in the diff).
$ws.Rows.Item(35).Insert()

$ws.Cells.Item(35, 1).Value = 6
$ws.Cells.Item(35, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(35, 3).Value = "Metropolitana"
$ws.Cells.Item(35, 4).Value = 44592
$ws.Cells.Item(35, 5).Value = 13
$ws.Cells.Item(35, 6).Value = 100112029
$ws.Cells.Item(35, 7).Value = "Orégano"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 33
$ws.Cells.Item(35, 11).Value = 8000
$ws.Cells.Item(35, 12).Value = 9000
$ws.Cells.Item(35, 13).Value = 8455
$ws.Cells.Item(35, 14).Value = "$/docena de atados"
$ws.Cells.Item(35, 15).Value = "Región Metropolitana"
$ws.Cells.Item(35, 16).Value = 2818
$ws.Cells.Item(35, 17).Value = 3
$ws.Cells.Item(35, 18).Value = "Hortaliza"
